$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing "İsim" data to column B
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("A1").Value = "Numara"

# New numeric values for rows 2-5, aligned with the names now in column B
$ws.Range("A2").Value = 201017
$ws.Range("A3").Value = 201014
$ws.Range("A4").Value = 201012
$ws.Range("A5").Value = 201015
